$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "GNG_TO-1651255489511969"
$wb.Worksheets.Item(2).Name = "NB_TO-16512554913136995"
$wb.Worksheets.Item(3).Name = "RS_TO-16512554913157086"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512554913777044"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512554914407024"

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651255489475971.csv"
$ws1.Range("B3").Value = "GNG_stims-16512554894949653.csv"
$ws1.Range("B4").Value = "go_stims-16512554894970007.csv"
$ws1.Range("B5").Value = "GNG_stims-1651255489510986.csv"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1651255490375018.csv"
$ws2.Range("B3").Value = "OB-16512554906497378.csv"
$ws2.Range("B4").Value = "ZB-match_1-1651255489520965.csv"
$ws2.Range("B5").Value = "TB-1651255490744738.csv"
$ws2.Range("B6").Value = "ZB-match_9-16512554895739727.csv"
$ws2.Range("B7").Value = "TB-16512554913017015.csv"
$ws2.Range("B8").Value = "TB-1651255490906705.csv"
$ws2.Range("B9").Value = "ZB-match_2-16512554896249664.csv"
$ws2.Range("B10").Value = "OB-16512554902530286.csv"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1651255491329739.csv"
$ws4.Range("B3").Value = "ZM_stims-16512554913177035.csv"
$ws4.Range("B4").Value = "MM_stims-16512554913607008.csv"
$ws4.Range("B5").Value = "ZM_stims-16512554913307035.csv"
$ws4.Range("B6").Value = "MM_stims-16512554913767078.csv"
$ws4.Range("B7").Value = "ZM_stims-16512554913617017.csv"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512554913817.csv"
$ws5.Range("B3").Value = "SAT_stims-16512554913927011.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512554914087362.csv"
$ws5.Range("B5").Value = "vSAT_stims-1651255491424735.csv"
